$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 2844.6667
$ws.Cells.Item(64, 9).Value = 2900
$ws.Cells.Item(64, 11).Value = 2900
$ws.Cells.Item(64, 13).Value = -2652

$ws.Cells.Item(67, 8).Value = 2844.6667
$ws.Cells.Item(67, 9).Value = 2900
$ws.Cells.Item(67, 11).Value = 2900
$ws.Cells.Item(67, 13).Value = -2042

$ws.Cells.Item(70, 8).Value = 1628.875
$ws.Cells.Item(70, 9).Value = 1241.8334
$ws.Cells.Item(70, 10).Value = 2790
$ws.Cells.Item(70, 11).Value = 3725.5002
$ws.Cells.Item(70, 12).Value = 8370
$ws.Cells.Item(70, 13).Value = -3455.5002
$ws.Cells.Item(70, 14).Value = -8910

$ws.Cells.Item(73, 8).Value = 1628.875
$ws.Cells.Item(73, 9).Value = 1241.8334
$ws.Cells.Item(73, 10).Value = 2790
$ws.Cells.Item(73, 11).Value = 3725.5002
$ws.Cells.Item(73, 12).Value = 8370
$ws.Cells.Item(73, 13).Value = -2789.5002
$ws.Cells.Item(73, 14).Value = -10242

$ws.Cells.Item(112, 8).Value = 1582.2941
$ws.Cells.Item(112, 10).Value = 1612.4375
$ws.Cells.Item(112, 12).Value = 4837.3125
$ws.Cells.Item(112, 14).Value = -7053.3125

$ws.Cells.Item(125, 8).Value = 142859000
$ws.Cells.Item(125, 9).Value = 500000500
$ws.Cells.Item(125, 10).Value = 2420
$ws.Cells.Item(125, 11).Value = 4500004500
$ws.Cells.Item(125, 12).Value = 21780
$ws.Cells.Item(125, 13).Value = -4500002040
$ws.Cells.Item(125, 14).Value = -26700

$ws.Cells.Item(127, 8).Value = 0
$ws.Cells.Item(127, 9).Value = 0
$ws.Cells.Item(127, 10).Value = 0
$ws.Cells.Item(127, 11).Value = 0
$ws.Cells.Item(127, 12).Value = 0
$ws.Cells.Item(127, 13).ClearContents()
$ws.Cells.Item(127, 14).ClearContents()

$ws.Cells.Item(128, 8).Value = 12500
$ws.Cells.Item(128, 10).Value = 12500
$ws.Cells.Item(128, 12).Value = 12500
$ws.Cells.Item(128, 14).Value = -22460

$ws.Cells.Item(137, 8).Value = 1307.8611
$ws.Cells.Item(137, 9).Value = 999.4138
$ws.Cells.Item(137, 10).Value = 2585.7144
$ws.Cells.Item(137, 11).Value = 2998.2414
$ws.Cells.Item(137, 12).Value = 7757.1432
$ws.Cells.Item(137, 13).Value = -448.2413999999999
$ws.Cells.Item(137, 14).Value = -12857.1432

$ws.Cells.Item(138, 8).Value = 2404.6775
$ws.Cells.Item(138, 10).Value = 3141.7322
$ws.Cells.Item(138, 12).Value = 9425.196599999999
$ws.Cells.Item(138, 14).Value = -19705.1966

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 746.0952
$ws.Cells.Item(74, 9).Value = 743.6129
$ws.Cells.Item(74, 10).Value = 900
$ws.Cells.Item(74, 11).Value = 743.6129
$ws.Cells.Item(74, 12).Value = 900
$ws.Cells.Item(74, 13).Value = 130.3871
$ws.Cells.Item(74, 14).Value = -2648

$ws.Cells.Item(77, 8).Value = 746.0952
$ws.Cells.Item(77, 9).Value = 743.6129
$ws.Cells.Item(77, 10).Value = 900
$ws.Cells.Item(77, 11).Value = 3718.0645
$ws.Cells.Item(77, 12).Value = 4500
$ws.Cells.Item(77, 13).Value = 649.9355
$ws.Cells.Item(77, 14).Value = -13236

$ws.Cells.Item(97, 8).Value = 2720
$ws.Cells.Item(97, 9).Value = 2720
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 11).Value = 2720
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 13).Value = -2224
$ws.Cells.Item(97, 14).ClearContents()

$ws.Cells.Item(122, 8).Value = 2042.862
$ws.Cells.Item(122, 9).Value = 1837.3182
$ws.Cells.Item(122, 10).Value = 2688.8572
$ws.Cells.Item(122, 11).Value = 5511.9546
$ws.Cells.Item(122, 12).Value = 8066.571599999999
$ws.Cells.Item(122, 13).Value = -3061.9546
$ws.Cells.Item(122, 14).Value = -12966.5716

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 2348.3333
$ws.Cells.Item(105, 9).Value = 2348.3333
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 11).Value = 2348.3333
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 13).Value = -601.3332999999998
$ws.Cells.Item(105, 14).ClearContents()

$ws.Cells.Item(107, 8).Value = 2145.875
$ws.Cells.Item(107, 9).Value = 2045.3077
$ws.Cells.Item(107, 10).Value = 2264.7273
$ws.Cells.Item(107, 11).Value = 2045.3077
$ws.Cells.Item(107, 12).Value = 2264.7273
$ws.Cells.Item(107, 13).Value = -125.3077000000001
$ws.Cells.Item(107, 14).Value = -6104.7273

$ws.Cells.Item(111, 8).Value = 29087.75
$ws.Cells.Item(111, 10).Value = 29087.75
$ws.Cells.Item(111, 12).Value = 29087.75
$ws.Cells.Item(111, 14).Value = -37267.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 894.1539
$ws.Cells.Item(16, 9).Value = 795.125
$ws.Cells.Item(16, 10).Value = 1052.6
$ws.Cells.Item(16, 11).Value = 795.125
$ws.Cells.Item(16, 12).Value = 1052.6
$ws.Cells.Item(16, 13).Value = -508.125
$ws.Cells.Item(16, 14).Value = -1626.6

$ws.Cells.Item(31, 8).Value = 5214.1377
$ws.Cells.Item(31, 9).Value = 4209.3335
$ws.Cells.Item(31, 10).Value = 6290.7144
$ws.Cells.Item(31, 11).Value = 4209.3335
$ws.Cells.Item(31, 12).Value = 6290.7144
$ws.Cells.Item(31, 13).Value = -3914.3335
$ws.Cells.Item(31, 14).Value = -6880.7144

$ws.Cells.Item(34, 8).Value = 5214.1377
$ws.Cells.Item(34, 9).Value = 4209.3335
$ws.Cells.Item(34, 10).Value = 6290.7144
$ws.Cells.Item(34, 11).Value = 4209.3335
$ws.Cells.Item(34, 12).Value = 6290.7144
$ws.Cells.Item(34, 13).Value = -4007.3335
$ws.Cells.Item(34, 14).Value = -6694.7144

$ws.Cells.Item(105, 8).Value = 551.6667
$ws.Cells.Item(105, 9).Value = 551.6667
$ws.Cells.Item(105, 11).Value = 551.6667
$ws.Cells.Item(105, 13).Value = 1195.3333

$ws.Cells.Item(113, 8).Value = 894.1539
$ws.Cells.Item(113, 9).Value = 795.125
$ws.Cells.Item(113, 10).Value = 1052.6
$ws.Cells.Item(113, 11).Value = 795.125
$ws.Cells.Item(113, 12).Value = 1052.6
$ws.Cells.Item(113, 13).Value = 1374.875
$ws.Cells.Item(113, 14).Value = -5392.6

$ws.Cells.Item(132, 8).Value = 4466874.5
$ws.Cells.Item(132, 9).Value = 1469.7858
$ws.Cells.Item(132, 10).Value = 8932279
$ws.Cells.Item(132, 11).Value = 4409.357400000001
$ws.Cells.Item(132, 12).Value = 26796837
$ws.Cells.Item(132, 13).Value = -1879.357400000001
$ws.Cells.Item(132, 14).Value = -26801897

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 449
$ws.Cells.Item(68, 9).Value = 449
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 11).Value = 1347
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 13).Value = -536
$ws.Cells.Item(68, 14).ClearContents()

$ws.Cells.Item(71, 8).Value = 449
$ws.Cells.Item(71, 9).Value = 449
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 11).Value = 4041
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 13).Value = 15
$ws.Cells.Item(71, 14).ClearContents()

$ws.Cells.Item(114, 8).Value = 5845.5
$ws.Cells.Item(114, 9).Value = 8692.923000000001
$ws.Cells.Item(114, 10).Value = 1732.5555
$ws.Cells.Item(114, 11).Value = 26078.769
$ws.Cells.Item(114, 12).Value = 5197.666499999999
$ws.Cells.Item(114, 13).Value = -22824.769
$ws.Cells.Item(114, 14).Value = -11705.6665

$ws.Cells.Item(131, 8).Value = 595.47
$ws.Cells.Item(131, 10).Value = 803.9032
$ws.Cells.Item(131, 12).Value = 2411.7096
$ws.Cells.Item(131, 14).Value = -12491.7096

$ws.Cells.Item(132, 8).Value = 1896.2963
$ws.Cells.Item(132, 10).Value = 2482.353
$ws.Cells.Item(132, 12).Value = 22341.177
$ws.Cells.Item(132, 14).Value = -27401.177

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1965.2667
$ws.Cells.Item(97, 9).Value = 2098.25
$ws.Cells.Item(97, 11).Value = 2098.25
$ws.Cells.Item(97, 13).Value = -1602.25

$ws.Cells.Item(112, 8).Value = 36124.125
$ws.Cells.Item(112, 10).Value = 36124.125
$ws.Cells.Item(112, 12).Value = 36124.125
$ws.Cells.Item(112, 14).Value = -38340.125

$ws.Cells.Item(126, 8).Value = 2508.3333
$ws.Cells.Item(126, 9).Value = 2400
$ws.Cells.Item(126, 10).Value = 2833.3333
$ws.Cells.Item(126, 11).Value = 7200
$ws.Cells.Item(126, 12).Value = 8499.999899999999
$ws.Cells.Item(126, 13).Value = -4730
$ws.Cells.Item(126, 14).Value = -13439.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 41669224
$ws.Cells.Item(7, 9).Value = 2500
$ws.Cells.Item(7, 10).Value = 62502588
$ws.Cells.Item(7, 11).Value = 2500
$ws.Cells.Item(7, 12).Value = 62502588
$ws.Cells.Item(7, 13).Value = -2388
$ws.Cells.Item(7, 14).Value = -62502812

$ws.Cells.Item(100, 8).Value = 333333340
$ws.Cells.Item(100, 9).Value = 333333340
$ws.Cells.Item(100, 11).Value = 333333340
$ws.Cells.Item(100, 13).Value = -333332799

$ws.Cells.Item(111, 8).Value = 29980
$ws.Cells.Item(111, 10).Value = 29980
$ws.Cells.Item(111, 12).Value = 29980
$ws.Cells.Item(111, 14).Value = -38160

$ws.Cells.Item(122, 8).Value = 8567.333000000001
$ws.Cells.Item(122, 9).Value = 11576
$ws.Cells.Item(122, 10).Value = 2550
$ws.Cells.Item(122, 11).Value = 34728
$ws.Cells.Item(122, 12).Value = 7650
$ws.Cells.Item(122, 13).Value = -32278
$ws.Cells.Item(122, 14).Value = -12550

$ws.Cells.Item(126, 8).Value = 41669224
$ws.Cells.Item(126, 9).Value = 2500
$ws.Cells.Item(126, 10).Value = 62502588
$ws.Cells.Item(126, 11).Value = 7500
$ws.Cells.Item(126, 12).Value = 187507764
$ws.Cells.Item(126, 13).Value = -5030
$ws.Cells.Item(126, 14).Value = -187512704

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(94, 8).Value = 18906
$ws.Cells.Item(94, 10).Value = 18906
$ws.Cells.Item(94, 12).Value = 18906
$ws.Cells.Item(94, 14).Value = -20708

$ws.Cells.Item(96, 8).Value = 125000400
$ws.Cells.Item(96, 9).Value = 125000400
$ws.Cells.Item(96, 11).Value = 125000400
$ws.Cells.Item(96, 13).Value = -124999027

$ws.Cells.Item(126, 8).Value = 35720820
$ws.Cells.Item(126, 9).Value = 52639308
$ws.Cells.Item(126, 10).Value = 4006.889
$ws.Cells.Item(126, 11).Value = 157917924
$ws.Cells.Item(126, 12).Value = 12020.667
$ws.Cells.Item(126, 13).Value = -157915454
$ws.Cells.Item(126, 14).Value = -16960.667
